# fix: prevent hidden columns from being labeled upon detecting changes (#11)
#
# Previously every comparison row in the merged AHB sheet got the "ANDERUNG"
# (change) marker in column L, even when the only differing cells were in a
# hidden column. The fix:
#   1. Clears the stale "ANDERUNG" label (and its bold/orange highlight
#      style) from column L wherever it was applied, resetting that cell to
#      the plain centered "no change" style.
#   2. For rows that are the first row of a new top-level segment group
#      (i.e. a genuine, visible difference was detected), re-applies the
#      shaded "group header" formatting across the whole row (A:V) so the
#      row is still highlighted as changed - matching the formatting
#      already used earlier in the sheet for the same kind of row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Template rows already carrying the desired formatting (rows 1-81 were
# already regenerated with the fix applied). Both templates live outside the
# 82-173 block being edited below so they stay stable as sources throughout.
#   - Row 2  : "group header" row style -> A:V = s2/s3(col B)/s4(col L)
#   - Row 2, col L : plain cleared/centered "no change" column-L style (s4)
$groupHeaderTemplate = $ws.Range("A2:V2")
$clearedLTemplate = $ws.Range("L2")

# Rows that are the first row of a new top-level segment group and need the
# whole row (A:V) re-shaded as "changed", with column L cleared.
$groupHeaderRows = @(82, 85, 89, 95, 100, 104, 108, 117, 123, 145, 168, 171)

foreach ($r in $groupHeaderRows) {
    $groupHeaderTemplate.Copy()
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial($xlPasteFormats)
    $ws.Range("L" + $r).ClearContents()
}

# Rows where only the column-L "ANDERUNG" marker needs to be removed because
# the detected change was actually confined to a hidden column; every other
# cell in the row keeps its existing formatting/content untouched.
$clearOnlyRowRanges = @(
    @(83, 84),
    @(86, 88),
    @(90, 94),
    @(96, 99),
    @(101, 103),
    @(105, 107),
    @(109, 116),
    @(118, 122),
    @(124, 144),
    @(146, 166),
    @(169, 170),
    @(172, 173)
)

foreach ($pair in $clearOnlyRowRanges) {
    $startRow = $pair[0]
    $endRow = $pair[1]
    $clearedLTemplate.Copy()
    $target = $ws.Range("L" + $startRow + ":L" + $endRow)
    $target.PasteSpecial($xlPasteFormats)
    $target.ClearContents()
}

$excel.CutCopyMode = 0
